$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 21

function Set-TextCell($col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# A: Code
Set-TextCell 1 "2033555832"
# B: Reference_Fournisseur (blank in source data)
Set-TextCell 2 ""
# C: Produits
Set-TextCell 3 "ddd"
# D: Unite_Stockage
Set-TextCell 4 "Unité"
# E: Unite_Commande
Set-TextCell 5 "Unité"
# F: Stock_Min
$ws.Cells.Item($row, 6).Value = 10
# G: Stock_Max
$ws.Cells.Item($row, 7).Value = 100
# H: Site
Set-TextCell 8 "Site principal"
# I: Lieu
Set-TextCell 9 "E2"
# J: Emplacement
Set-TextCell 10 "E2"
# K: Fournisseur
Set-TextCell 11 "FournX"
# L: Prix_Unitaire
$ws.Cells.Item($row, 12).Value = 0
# M: Categorie
Set-TextCell 13 "Général"
# N: Secteur
Set-TextCell 14 "Général"
# O: Reference
Set-TextCell 15 "2033555832"
# P: Quantite
$ws.Cells.Item($row, 16).Value = 0
# Q: Date_Entree
Set-TextCell 17 "2025-06-04"
